{"js": "const body = context.document.body;\n\n// The placeholder paragraph is the very first paragraph in the document.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nconst firstParagraph = paragraphs.items[0];\n\n// Rebuild the paragraph: add a paragraph border (top/left/bottom/right, 5pt\n// spacing, no visible line), bump the left indent from 120 -> 225 twips, and\n// collapse the two runs (\"**ID__...topic_13__ID**\" + a trailing space run)\n// into a single run carrying the new placeholder id with no trailing space.\nconst newParagraphOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n      'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' +\n            '<w:p>' +\n              '<w:pPr>' +\n                '<w:pBdr>' +\n                  '<w:top w:space=\"5\"/>' +\n                  '<w:left w:space=\"5\"/>' +\n                  '<w:bottom w:space=\"5\"/>' +\n                  '<w:right w:space=\"5\"/>' +\n                '</w:pBdr>' +\n                '<w:spacing w:after=\"0\"/>' +\n                '<w:ind w:left=\"225\"/>' +\n                '<w:jc w:val=\"left\"/>' +\n              '</w:pPr>' +\n              '<w:r>' +\n                '<w:rPr>' +\n                  '<w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\"/>' +\n                  '<w:b w:val=\"false\"/>' +\n                  '<w:i w:val=\"false\"/>' +\n                  '<w:color w:val=\"000000\"/>' +\n                  '<w:sz w:val=\"22\"/>' +\n                '</w:rPr>' +\n                '<w:t>**ID__AFFARS_MP_5315_3_3_1__ID**</w:t>' +\n              '</w:r>' +\n            '</w:p>' +\n          '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\nfirstParagraph.insertOoxml(newParagraphOoxml, \"Replace\");\nawait context.sync();\n\n// insertOoxml drops explicit \"false\" boolean run properties (they come back\n// out as simply absent, which already renders as not-bold/not-italic); push\n// italic back onto the run explicitly so the rPr keeps matching the source\n// formatting as closely as the API allows.\nconst newRuns = body.search(\"**ID__AFFARS_MP_5315_3_3_1__ID**\", { matchCase: true });\nnewRuns.load(\"items\");\nawait context.sync();\nnewRuns.items[0].font.set({ italic: false });\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# First paragraph of the document: the \"**ID__AFFARS_...\" placeholder line.\n$para = $d.Paragraphs(1)\n\n# Replace the paragraph's text (both runs combine to \"**ID__...topic_13__ID** \")\n# with the new id, dropping the trailing space run entirely.\n$r = $para.Range\n$r.End = $r.End - 1\n$r.Text = \"**ID__AFFARS_MP_5315_3_3_1__ID**\"\n\n# Add a paragraph border (top/left/bottom/right) with 5pt-equivalent spacing\n# and no visible line - matches <w:pBdr><w:top w:space=\"5\"/>...</w:pBdr>.\n$para.Range.ParagraphFormat.Borders.DistanceFromTop = 5\n$para.Range.ParagraphFormat.Borders.DistanceFromLeft = 5\n$para.Range.ParagraphFormat.Borders.DistanceFromBottom = 5\n$para.Range.ParagraphFormat.Borders.DistanceFromRight = 5\n\n# Change the left indent from 120 twips (6pt) to 225 twips (11.25pt).\n$para.Range.ParagraphFormat.LeftIndent = 11.25\n"}
